# Update gh-pages output data (view counts / attendee figures) on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.
#
# Cell  Old   New
# F6    98    100
# F10   4913  4925
# F11   4618  4627
# F13   27    28
# F15   37    38

$wb = $excel.ActiveWorkbook

$updates = @{
    "F6"  = 100
    "F10" = 4925
    "F11" = 4627
    "F13" = 28
    "F15" = 38
}

$targetSheets = @(1, 4)

foreach ($idx in $targetSheets) {
    $ws = $wb.Worksheets.Item($idx)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
